$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.287.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.238.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.10%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.75%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -2.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.237.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.129'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.08%  '

$ws.Range("E11").Value = '  +3.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.411'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.802.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.10%  '

$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.78'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.221.78'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000167'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.242.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.26%  '

$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.507'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.380.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000118'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.181'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.48%  '

$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("E31").Value = '  +5.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.62%  '

$ws.Range("E34").Value = '  +0.05%  '

$ws.Range("E35").Value = '  +7.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.79'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.84'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.48'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.860'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.62%  '

$ws.Range("E40").Value = '  +10.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +14.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '362.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +16.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.710.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.19%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.92%  '

$ws.Range("E49").Value = '  +5.11%  '

$ws.Range("E50").Value = '  +3.43%  '

$ws.Range("E51").Value = '  +1.66%  '
